$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1493, 114, 7, 18),
    @(1493, 114, 7, 18),
    @(1493, 114, 7, 18),
    @(1493, 114, 7, 18),
    @(1493, 114, 7, 18),
    @(1493, 114, 7, 18),
    @(1493, 114, 7, 18),
    @(1493, 114, 7, 18),
    @(1493, 114, 7, 20),
    @(2000, 149, 7, 20)
)

$row = 3
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row = $row + 1
}
